# Apply the "updated task used in testing" edit:
#   - D2: 6 -> 7
#   - F2: 2 -> 3
#   - H2: 36 -> 46
#   - active selection moves from D5 to D2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 7
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 46

# Update the selected/active cell shown in the saved sheet view.
$ws.Range("D2").Select()
